# Applies the cryptos list update (prices / 1h volume deltas) as of
# Thu Jan  4 23:08:55 UTC 2024, generated from the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "321.91") must be
# forced to Text format first, otherwise Excel will silently convert the
# assigned string into a floating point number (losing the exact original
# formatting, e.g. trailing zeros). We flip the format to Text, write the
# value, then restore the default "Normal" style so no stray formatting is
# left behind on the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '44.446.67'
$ws.Range('E2').Value = '  +3.82%  '
$ws.Range('D3').Value = '2.270.96'
$ws.Range('E3').Value = '  +3.01%  '
$ws.Range('E4').Value = '  -0.12%  '
Set-TextValue 'D5' '321.91'
$ws.Range('E5').Value = '  +2.00%  '
Set-TextValue 'D6' '104.88'
$ws.Range('E6').Value = '  +6.03%  '
Set-TextValue 'D7' '0.590'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  -0.05%  '
Set-TextValue 'D9' '0.571'
$ws.Range('E9').Value = '  +2.03%  '
Set-TextValue 'D10' '38.63'
$ws.Range('E10').Value = '  +5.25%  '
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('E13').Value = '  +0.57%  '
Set-TextValue 'D14' '0.884'
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('D15').Value = '2.617.77'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').Value = '2.265.76'
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').Value = '44.343.03'
$ws.Range('E18').Value = '  +3.76%  '
Set-TextValue 'D19' '13.99'
$ws.Range('E19').Value = '  -3.42%  '
$ws.Range('E20').Value = '  +4.74%  '
Set-TextValue 'D21' '6.54'
$ws.Range('E21').Value = '  +2.12%  '
Set-TextValue 'D22' '66.36'
$ws.Range('E22').Value = '  +1.87%  '
Set-TextValue 'D23' '3.20'
$ws.Range('E23').Value = '  +1.60%  '
Set-TextValue 'D24' '239.65'
$ws.Range('E24').Value = '  +1.64%  '
Set-TextValue 'D25' '2.22'
$ws.Range('E25').Value = '  +4.51%  '
$ws.Range('E26').Value = '  -0.01%  '
Set-TextValue 'D27' '10.20'
$ws.Range('E27').Value = '  +2.34%  '
Set-TextValue 'D28' '38.59'
$ws.Range('E28').Value = '  +13.15%  '
Set-TextValue 'D29' '2.21'
$ws.Range('E29').Value = '  -0.32%  '
Set-TextValue 'D30' '6.49'
$ws.Range('E30').Value = '  +3.45%  '
Set-TextValue 'D31' '20.66'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D32' '0.0884'
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D33' '161.63'
$ws.Range('E33').Value = '  +3.81%  '
Set-TextValue 'D34' '2.77'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E35').Value = '  +9.59%  '
$ws.Range('E36').Value = '  +5.43%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  +0.49%  '
Set-TextValue 'D39' '3.94'
$ws.Range('E39').Value = '  +3.16%  '
Set-TextValue 'D40' '4.44'
$ws.Range('E40').Value = '  +0.60%  '
Set-TextValue 'D41' '15.61'
$ws.Range('E41').Value = '  +26.63%  '
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '1.775.18'
$ws.Range('E44').Value = '  -5.46%  '
$ws.Range('E45').Value = '  +0.88%  '
Set-TextValue 'D46' '86.58'
$ws.Range('E46').Value = '  -1.07%  '
Set-TextValue 'D47' '5.45'
$ws.Range('E47').Value = '  +2.13%  '
Set-TextValue 'D48' '60.63'
$ws.Range('E48').Value = '  -0.20%  '
Set-TextValue 'D49' '75.07'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('E50').Value = '  +7.62%  '
Set-TextValue 'D51' '104.10'
$ws.Range('E51').Value = '  +2.06%  '

Write-Host "Applied cryptos list update (88 cell changes across rows 2-51)"
